# Completed the vote validation module: refresh class_summary_df_yolo
# validation metrics (TP/FP/FN/Total/Precision/Recall per class).
# The "stamp" class was split into "invalid_stamp" (inserted before "key")
# and "valid_stamp" (inserted before "water_glass"), and a new class
# "wooden_wheel" row was appended (A1:H44 -> A1:H45), with every row's
# TP/FP/FN/Total/Precision/Recall recomputed against the new validation run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2: balance
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "balance"
$ws.Range("C2").Value = 14
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 16
$ws.Range("G2").Value = 0.875
$ws.Range("H2").Value = 1

# row 3: bus
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "bus"
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 18
$ws.Range("G3").Value = 0.7647058823529411
$ws.Range("H3").Value = 0.9285714285714286

# row 4: candle_light
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "candle_light"
$ws.Range("C4").Value = 14
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 14
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1

# row 5: computer
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "computer"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 14
$ws.Range("G5").Value = 0.7142857142857143
$ws.Range("H5").Value = 0.4166666666666667

# row 6: conch
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "conch"
$ws.Range("C6").Value = 13
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 14
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 0.9285714285714286

# row 7: cycle
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "cycle"
$ws.Range("C7").Value = 14
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 14
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1

# row 8: damphu
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "damphu"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 14
$ws.Range("F8").Value = 14
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = 0

# row 9: dog
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "dog"
$ws.Range("C9").Value = 14
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 23
$ws.Range("G9").Value = 0.6086956521739131
$ws.Range("H9").Value = 1

# row 10: farmer
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "farmer"
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 14
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0.2857142857142857

# row 11: hammer_scythe
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "hammer_scythe"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 14
$ws.Range("F11").Value = 14
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = 0

# row 12: hand_palm
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "hand_palm"
$ws.Range("C12").Value = 14
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 14
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1

# row 13: heart
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "heart"
$ws.Range("C13").Value = 14
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 14
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 1

# row 14: hoe
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "hoe"
$ws.Range("C14").Value = 13
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 15
$ws.Range("G14").Value = 0.9285714285714286
$ws.Range("H14").Value = 0.9285714285714286

# row 15: house
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "house"
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 14
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0.3571428571428572

# row 16: invalid_stamp
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "invalid_stamp"
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 0.6
$ws.Range("H16").Value = 0.5

# row 17: key
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "key"
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1

# row 18: ladder
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "ladder"
$ws.Range("C18").Value = 14
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 1

# row 19: lock
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "lock"
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 14
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 0.9285714285714286

# row 20: lotus
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "lotus"
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 0.7
$ws.Range("H20").Value = 1

# row 21: loud_speaker
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "loud_speaker"
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 14
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 1

# row 22: mother_and_child
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "mother_and_child"
$ws.Range("C22").Value = 11
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 14
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0.7857142857142857

# row 23: namaste
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "namaste"
$ws.Range("C23").Value = 14
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 1

# row 24: nepali_big_basket
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "nepali_big_basket"
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 14
$ws.Range("F24").Value = 14
$ws.Range("G24").Value = ""
$ws.Range("H24").Value = 0

# row 25: nepali_cap
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "nepali_cap"
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 6
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 0.8888888888888888
$ws.Range("H25").Value = 0.5714285714285714

# row 26: nepali_jug
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "nepali_jug"
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 9
$ws.Range("F26").Value = 14
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0.3571428571428572

# row 27: nepali_madal
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "nepali_madal"
$ws.Range("C27").Value = 9
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = 16
$ws.Range("G27").Value = 0.8181818181818182
$ws.Range("H27").Value = 0.6428571428571429

# row 28: nepali_small_basket
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "nepali_small_basket"
$ws.Range("C28").Value = 14
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 14
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 1

# row 29: owl
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "owl"
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 14
$ws.Range("F29").Value = 14
$ws.Range("G29").Value = ""
$ws.Range("H29").Value = 0

# row 30: pen
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "pen"
$ws.Range("C30").Value = 12
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 2
$ws.Range("F30").Value = 14
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0.8571428571428571

# row 31: roaster
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "roaster"
$ws.Range("C31").Value = 14
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 15
$ws.Range("G31").Value = 0.9333333333333333
$ws.Range("H31").Value = 1

# row 32: sheep
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "sheep"
$ws.Range("C32").Value = 11
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 3
$ws.Range("F32").Value = 14
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 0.7857142857142857

# row 33: star
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = "star"
$ws.Range("C33").Value = 14
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 14
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 1

# row 34: stick
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = "stick"
$ws.Range("C34").Value = 7
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 7
$ws.Range("F34").Value = 14
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 0.5

# row 35: sun
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = "sun"
$ws.Range("C35").Value = 14
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 14
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 1

# row 36: tiger
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = "tiger"
$ws.Range("C36").Value = 12
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 1
$ws.Range("F36").Value = 14
$ws.Range("G36").Value = 0.9230769230769231
$ws.Range("H36").Value = 0.9230769230769231

# row 37: torch_light
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = "torch_light"
$ws.Range("C37").Value = 9
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 4
$ws.Range("F37").Value = 14
$ws.Range("G37").Value = 0.9
$ws.Range("H37").Value = 0.6923076923076923

# row 38: tree
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = "tree"
$ws.Range("C38").Value = 11
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 3
$ws.Range("F38").Value = 14
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 0.7857142857142857

# row 39: turtle
$ws.Range("A39").Value = 37
$ws.Range("B39").Value = "turtle"
$ws.Range("C39").Value = 5
$ws.Range("D39").Value = 2
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = 16
$ws.Range("G39").Value = 0.7142857142857143
$ws.Range("H39").Value = 0.3571428571428572

# row 40: umbrella
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = "umbrella"
$ws.Range("C40").Value = 14
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 14
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 1

# row 41: valid_stamp
$ws.Range("A41").Value = 39
$ws.Range("B41").Value = "valid_stamp"
$ws.Range("C41").Value = 4
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 4
$ws.Range("F41").Value = 9
$ws.Range("G41").Value = 0.8
$ws.Range("H41").Value = 0.5

# row 42: water_glass
$ws.Range("A42").Value = 40
$ws.Range("B42").Value = "water_glass"
$ws.Range("C42").Value = 14
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 14
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 1

# row 43: water_jug
$ws.Range("A43").Value = 41
$ws.Range("B43").Value = "water_jug"
$ws.Range("C43").Value = 11
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 3
$ws.Range("F43").Value = 14
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 0.7857142857142857

# row 44: woman_man
$ws.Range("A44").Value = 42
$ws.Range("B44").Value = "woman_man"
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 14
$ws.Range("F44").Value = 14
$ws.Range("G44").Value = ""
$ws.Range("H44").Value = 0

# row 45: wooden_wheel
$ws.Range("A44").Copy() | Out-Null
$ws.Range("A45").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = "wooden_wheel"
$ws.Range("C45").Value = 14
$ws.Range("D45").Value = 5
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 19
$ws.Range("G45").Value = 0.7368421052631579
$ws.Range("H45").Value = 1
